$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 136, shifting old rows 136-137
# down to 137-138 (their contents stay intact thanks to the native
# row-insert semantics).
$ws.Rows.Item(136).Insert()

# Populate the freshly inserted row 136 with the new weekly record.
$ws.Cells.Item(136, 1).Value = 10
$ws.Cells.Item(136, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(136, 3).Value = "La Araucanía"
$ws.Cells.Item(136, 4).Value = 44448
$ws.Cells.Item(136, 5).Value = 9
$ws.Cells.Item(136, 6).Value = 100112052
$ws.Cells.Item(136, 7).Value = "Albahaca"
$ws.Cells.Item(136, 8).Value = "Sin especificar"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 40
$ws.Cells.Item(136, 11).Value = 7000
$ws.Cells.Item(136, 12).Value = 7000
$ws.Cells.Item(136, 13).Value = 7000
$ws.Cells.Item(136, 14).Value = "`$/paquete"
$ws.Cells.Item(136, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(136, 16).Value = 7000
$ws.Cells.Item(136, 17).Value = 1
$ws.Cells.Item(136, 18).Value = "Hortaliza"
